$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "expected result" text for step 1 (G3) to the new, expanded
# two-line wording that covers both the pop-up appearing and not appearing.
$newText = "1. Страница https://market.yandex.ru/ открыта, появился поп-ап с предложением авторизации`n2. Страница https://market.yandex.ru/ открыта, поп-ап с предложением авторизации не появился - перейти к шагу 3.`n"
$ws.Range("G3").Value = $newText

# The extra line makes row 3 taller, matching the height used by row 6
# (which also wraps two lines of text at this column width).
$ws.Rows.Item(3).RowHeight = 102

# Selection/view state as left by the editor after making the change.
$ws.Range("H3").Select()
